$d = $word.ActiveDocument

function Insert-PastTenseD($searchPhrase) {
    $r = $d.Content
    $r.Find.Execute($searchPhrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($r.Find.Found) {
        $insertPos = $r.End
        $insertPoint = $d.Range($insertPos, $insertPos)
        $insertPoint.InsertAfter("d")
        # Toggle formatting on just the inserted "d" so it stays its own run
        # (identical rPr to its neighbors) instead of being merged back in.
        $dRange = $d.Range($insertPos, $insertPos + 1)
        $dRange.Font.Bold = 1
        $dRange.Font.Bold = 0
    }
}

# "Primary duties include response..." -> "Primary duties included response..."
Insert-PastTenseD("Primary duties include")

# "...that each include replication..." -> "...that each included replication..."
Insert-PastTenseD("that each include")
